$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (J1:M1 renamed) ---
$ws.Range("J1").Value = "CET1_resource"
$ws.Range("K1").Value = "T1_resource"
$ws.Range("L1").Value = "total_capital_resource"
$ws.Range("M1").Value = "TLAC_resource"

# --- Prepare style template for new index-column rows (copy format of A2) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A28").PasteSpecial(-4122) | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null

# --- Row 2 (prime_auto) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "prime_auto"
$ws.Range("C2").Value = 15000
$ws.Range("D2").Value = 8000
$ws.Range("E2").Value = -5000
$ws.Range("F2").Value = 80
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.2
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 23000
$ws.Range("O2").Value = 10000
$ws.Range("P2").Value = "(10000, 23000)"
$ws.Range("Q2").Value = 10000
$ws.Range("R2").Value = -5000

# --- Row 3 (subprime_auto) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "subprime_auto"
$ws.Range("C3").Value = 4000
$ws.Range("D3").Value = 2000
$ws.Range("E3").Value = -1000
$ws.Range("F3").Value = 150
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0.6
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = "(3000, 6000)"
$ws.Range("Q3").Value = 6000
$ws.Range("R3").Value = 2000

# --- Row 4 (mtg_30_fixed) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "mtg_30_fixed"
$ws.Range("C4").Value = 50000
$ws.Range("D4").Value = 30000
$ws.Range("E4").Value = -15000
$ws.Range("F4").Value = 70
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0.3
$ws.Range("I4").Value = 0.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 80000
$ws.Range("O4").Value = 35000
$ws.Range("P4").Value = "(35000, 80000)"
$ws.Range("Q4").Value = 80000
$ws.Range("R4").Value = 30000

# --- Row 5 (mtg_15_fixed) ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "mtg_15_fixed"
$ws.Range("C5").Value = 15000
$ws.Range("D5").Value = 8000
$ws.Range("E5").Value = -3000
$ws.Range("F5").Value = 72
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0.3
$ws.Range("I5").Value = 0.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 23000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = "(12000, 23000)"
$ws.Range("Q5").Value = 23000
$ws.Range("R5").Value = 8000

# --- Row 6 (mtg_7_fixed) ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "mtg_7_fixed"
$ws.Range("C6").Value = 5000
$ws.Range("D6").Value = 2500
$ws.Range("E6").Value = -2000
$ws.Range("F6").Value = 70
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0.3
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 7500
$ws.Range("O6").Value = 3000
$ws.Range("P6").Value = "(3000, 7500)"
$ws.Range("Q6").Value = 7500
$ws.Range("R6").Value = 2500

# --- Row 7 (mtg_15_arm) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "mtg_15_arm"
$ws.Range("C7").Value = 20000
$ws.Range("D7").Value = 10000
$ws.Range("E7").Value = -4000
$ws.Range("F7").Value = 71
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0.3
$ws.Range("I7").Value = 0.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 30000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = "(16000, 30000)"
$ws.Range("Q7").Value = 30000
$ws.Range("R7").Value = 10000

# --- Row 8 (mtg_7_arm) ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "mtg_7_arm"
$ws.Range("C8").Value = 10000
$ws.Range("D8").Value = 3500
$ws.Range("E8").Value = -2000
$ws.Range("F8").Value = 75
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0.3
$ws.Range("I8").Value = 0.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 13500
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = "(8000, 13500)"
$ws.Range("Q8").Value = 13500
$ws.Range("R8").Value = 3500

# --- Row 9 (consumer_card) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "consumer_card"
$ws.Range("C9").Value = 60000
$ws.Range("D9").Value = 20000
$ws.Range("E9").Value = -10000
$ws.Range("F9").Value = 550
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1.1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 80000
$ws.Range("O9").Value = 50000
$ws.Range("P9").Value = "(50000, 80000)"
$ws.Range("Q9").Value = 80000
$ws.Range("R9").Value = 20000

# --- Row 10 (business_card) ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "business_card"
$ws.Range("C10").Value = 15000
$ws.Range("D10").Value = 2500
$ws.Range("E10").Value = -2000
$ws.Range("F10").Value = 400
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0.95
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 17500
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = "(13000, 17500)"
$ws.Range("Q10").Value = 17500
$ws.Range("R10").Value = 2500

# --- Row 11 (business_loan_revolver) ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "business_loan_revolver"
$ws.Range("C11").Value = 20000
$ws.Range("D11").Value = 15000
$ws.Range("E11").Value = -5000
$ws.Range("F11").Value = 110
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 0.75
$ws.Range("I11").Value = 1.2
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 35000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = "(15000, 35000)"
$ws.Range("Q11").Value = 15000
$ws.Range("R11").Value = -5000

# --- Row 12 (business_loan_term) ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "business_loan_term"
$ws.Range("C12").Value = 55000
$ws.Range("D12").Value = 15000
$ws.Range("E12").Value = -8500
$ws.Range("F12").Value = 100
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0.7
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 70000
$ws.Range("O12").Value = 46500
$ws.Range("P12").Value = "(46500, 70000)"
$ws.Range("Q12").Value = 51562.09677419356
$ws.Range("R12").Value = -3437.90322580644

# --- Row 13 (commercial_loan_revolver) ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "commercial_loan_revolver"
$ws.Range("C13").Value = 40000
$ws.Range("D13").Value = 20000
$ws.Range("E13").Value = -12000
$ws.Range("F13").Value = 95
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 0.65
$ws.Range("I13").Value = 1.15
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 60000
$ws.Range("O13").Value = 28000
$ws.Range("P13").Value = "(28000, 60000)"
$ws.Range("Q13").Value = 28000
$ws.Range("R13").Value = -12000

# --- Row 14 (commercial_loan_term) ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "commercial_loan_term"
$ws.Range("C14").Value = 100000
$ws.Range("D14").Value = 40000
$ws.Range("E14").Value = -15000
$ws.Range("F14").Value = 85
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0.6
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 140000
$ws.Range("O14").Value = 85000
$ws.Range("P14").Value = "(85000, 140000)"
$ws.Range("Q14").Value = 85000
$ws.Range("R14").Value = -15000

# --- Row 15 (rates) ---
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "rates"
$ws.Range("C15").Value = 80000
$ws.Range("D15").Value = 25000
$ws.Range("E15").Value = -10000
$ws.Range("F15").Value = 112
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0.4
$ws.Range("I15").Value = 0.55
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 105000
$ws.Range("O15").Value = 70000
$ws.Range("P15").Value = "(70000, 105000)"
$ws.Range("Q15").Value = 105000
$ws.Range("R15").Value = 25000

# --- Row 16 (equities) ---
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "equities"
$ws.Range("C16").Value = 60000
$ws.Range("D16").Value = 15000
$ws.Range("E16").Value = -12000
$ws.Range("F16").Value = 108
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 0.5
$ws.Range("I16").Value = 0.6
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 75000
$ws.Range("O16").Value = 48000
$ws.Range("P16").Value = "(48000, 75000)"
$ws.Range("Q16").Value = 75000
$ws.Range("R16").Value = 15000

# --- Row 17 (fixed_income) ---
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "fixed_income"
$ws.Range("C17").Value = 75000
$ws.Range("D17").Value = 20000
$ws.Range("E17").Value = -10000
$ws.Range("F17").Value = 140
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 0.35
$ws.Range("I17").Value = 0.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 95000
$ws.Range("O17").Value = 65000
$ws.Range("P17").Value = "(65000, 95000)"
$ws.Range("Q17").Value = 95000
$ws.Range("R17").Value = 20000

# --- Row 18 (prime) ---
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "prime"
$ws.Range("C18").Value = 30000
$ws.Range("D18").Value = 4000
$ws.Range("E18").Value = -2500
$ws.Range("F18").Value = 160
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 0.5
$ws.Range("I18").Value = 0.7
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 34000
$ws.Range("O18").Value = 27500
$ws.Range("P18").Value = "(27500, 34000)"
$ws.Range("Q18").Value = 34000
$ws.Range("R18").Value = 4000

# --- Row 19 (payments) ---
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "payments"
$ws.Range("C19").Value = 1000
$ws.Range("D19").Value = 500
$ws.Range("E19").Value = -250
$ws.Range("F19").Value = 600
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 0.05
$ws.Range("I19").Value = 0.05
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 1500
$ws.Range("O19").Value = 750
$ws.Range("P19").Value = "(750, 1500)"
$ws.Range("Q19").Value = 1500
$ws.Range("R19").Value = 500

# --- Row 20 (consumer_checking) ---
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "consumer_checking"
$ws.Range("C20").Value = 175000
$ws.Range("D20").Value = 60000
$ws.Range("E20").Value = -40000
$ws.Range("F20").Value = 200
$ws.Range("G20").Value = -1
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 235000
$ws.Range("O20").Value = 135000
$ws.Range("P20").Value = "(135000, 235000)"
$ws.Range("Q20").Value = 235000
$ws.Range("R20").Value = 60000

# --- Row 21 (consumer_savings) ---
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "consumer_savings"
$ws.Range("C21").Value = 60000
$ws.Range("D21").Value = 25000
$ws.Range("E21").Value = -10000
$ws.Range("F21").Value = 185
$ws.Range("G21").Value = -1
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 85000
$ws.Range("O21").Value = 50000
$ws.Range("P21").Value = "(50000, 85000)"
$ws.Range("Q21").Value = 85000
$ws.Range("R21").Value = 25000

# --- Row 22 (business_op_deposit) ---
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "business_op_deposit"
$ws.Range("C22").Value = 100000
$ws.Range("D22").Value = 30000
$ws.Range("E22").Value = -15000
$ws.Range("F22").Value = 150
$ws.Range("G22").Value = -1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 130000
$ws.Range("O22").Value = 85000
$ws.Range("P22").Value = "(85000, 130000)"
$ws.Range("Q22").Value = 130000
$ws.Range("R22").Value = 30000

# --- Row 23 (business_nonop_deposit) ---
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "business_nonop_deposit"
$ws.Range("C23").Value = 20000
$ws.Range("D23").Value = 15000
$ws.Range("E23").Value = -3000
$ws.Range("F23").Value = 25
$ws.Range("G23").Value = -1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 35000
$ws.Range("O23").Value = 17000
$ws.Range("P23").Value = "(17000, 35000)"
$ws.Range("Q23").Value = 17000
$ws.Range("R23").Value = -3000

# --- Row 24 (commercial_op_deposits) ---
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "commercial_op_deposits"
$ws.Range("C24").Value = 90500
$ws.Range("D24").Value = 35000
$ws.Range("E24").Value = -10000
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = -1
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 125500
$ws.Range("O24").Value = 80500
$ws.Range("P24").Value = "(80500, 125500)"
$ws.Range("Q24").Value = 125500
$ws.Range("R24").Value = 35000

# --- Row 25 (commercial_nonop_deposit) ---
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "commercial_nonop_deposit"
$ws.Range("C25").Value = 42000
$ws.Range("D25").Value = 20000
$ws.Range("E25").Value = -4000
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = -1
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 62000
$ws.Range("O25").Value = 38000
$ws.Range("P25").Value = "(38000, 62000)"
$ws.Range("Q25").Value = 38000
$ws.Range("R25").Value = -4000

# --- Row 26 (commercial_paper) ---
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "commercial_paper"
$ws.Range("C26").Value = 53000
$ws.Range("D26").Value = 1000000
$ws.Range("E26").Value = -50000
$ws.Range("F26").Value = -20
$ws.Range("G26").Value = -1
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 1053000
$ws.Range("O26").Value = 3000
$ws.Range("P26").Value = "(3000, 1053000)"
$ws.Range("Q26").Value = 3000
$ws.Range("R26").Value = -50000

# --- Row 27 (equity) ---
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "equity"
$ws.Range("C27").Value = 52500
$ws.Range("D27").Value = 100000
$ws.Range("E27").Value = -52500
$ws.Range("F27").Value = -1000
$ws.Range("G27").Value = -1
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1
$ws.Range("K27").Value = 1
$ws.Range("L27").Value = 1
$ws.Range("M27").Value = 1
$ws.Range("N27").Value = 152500
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = "(0, 152500)"
$ws.Range("Q27").Value = 60652.58064516129
$ws.Range("R27").Value = 8152.580645161288

# --- Row 28 (prefs) ---
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "prefs"
$ws.Range("C28").Value = 8000
$ws.Range("D28").Value = 13000
$ws.Range("E28").Value = -7000
$ws.Range("F28").Value = -650
$ws.Range("G28").Value = -1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 1
$ws.Range("L28").Value = 1
$ws.Range("M28").Value = 1
$ws.Range("N28").Value = 21000
$ws.Range("O28").Value = 1000
$ws.Range("P28").Value = "(1000, 21000)"
$ws.Range("Q28").Value = 8270.806451612903
$ws.Range("R28").Value = 270.8064516129034

# --- Row 29 (sub_debt) ---
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "sub_debt"
$ws.Range("C29").Value = 11000
$ws.Range("D29").Value = 20000
$ws.Range("E29").Value = -9500
$ws.Range("F29").Value = -200
$ws.Range("G29").Value = -1
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 1
$ws.Range("M29").Value = 1
$ws.Range("N29").Value = 31000
$ws.Range("O29").Value = 1500
$ws.Range("P29").Value = "(1500, 31000)"
$ws.Range("Q29").Value = 11027.74193548387
$ws.Range("R29").Value = 27.7419354838712

# --- Row 30 (senior_debt) ---
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = "senior_debt"
$ws.Range("C30").Value = 43000
$ws.Range("D30").Value = 60000
$ws.Range("E30").Value = -40000
$ws.Range("F30").Value = -120
$ws.Range("G30").Value = -1
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 1
$ws.Range("N30").Value = 103000
$ws.Range("O30").Value = 3000
$ws.Range("P30").Value = "(3000, 103000)"
$ws.Range("Q30").Value = 44110.96774193548
$ws.Range("R30").Value = 1110.967741935485
